$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.898.88"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "3.435.25"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.41"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.34"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").Value = "3.436.09"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.77"
$ws.Range("E10").Value = "  +1.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("E12").Value = "  +2.82%  "

$ws.Range("D13").Value = "4.024.62"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("E14").Value = "  +2.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.92"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "3.427.63"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("D18").Value = "62.928.08"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.38"
$ws.Range("E19").Value = "  +2.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.18"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.60"
$ws.Range("E22").Value = "  -3.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.560"
$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.43"
$ws.Range("E24").Value = "  -1.40%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").Value = "3.594.25"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("E27").Value = "  -3.92%  "

$ws.Range("E28").Value = "  -5.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.55"
$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.10"
$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.24"
$ws.Range("E34").Value = "  -2.57%  "

$ws.Range("E35").Value = "  -7.99%  "

$ws.Range("E36").Value = "  -1.01%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  -1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.63"
$ws.Range("E39").Value = "  +4.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "168.14"
$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").Value = "3.472.95"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.787"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.34"
$ws.Range("E44").Value = "  -1.33%  "

$ws.Range("E45").Value = "  -0.83%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("E47").Value = "  -3.09%  "

$ws.Range("D48").Value = "2.574.26"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("E49").Value = "  +3.07%  "

$ws.Range("E50").Value = "  +0.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.63"
$ws.Range("E51").Value = "  -3.75%  "
